$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "release/1.0.4"
$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("D6").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("F6").Value = "X"
